$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 24: "875. Koko Eating Bananas" (Binary Search category) ---
$question = "875. Koko Eating Bananas"
$solution = @'
Do binary search ON "k" dont sort any piles or anything
kstart =1, kend = max(piles)
use while loop (kstart<=kend){
kmid = middle of kstart and kend
calculate timeTaken with kmid
if(timeTake<=h) { //we are taking less than expected time, so we can further decrease speed of eating
	k = min(k, kmid)
	kend = kmid -1
}
else { //we are taking more than expected time, so increase speed of eating
k start = kmid+1
}
'@

# Copy the formatting of the previous data row (row 23) down onto the new
# row 24 so fills/borders/alignment match the rest of the table.
$ws.Range("A23:C23").Copy()
$ws.Range("A24:C24").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("A24").Value = "Binary Search"
$ws.Range("B24").Value = $question
$ws.Range("C24").Value = $solution

# The wrapped, multi-line solution text needs a taller row (Excel would
# auto-calculate this from the wrapped line count); set it explicitly.
$ws.Rows(24).RowHeight = 201.6

# Match the author's final selection / active cell.
$ws.Range("C24").Select()
